$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch two scratch cells to register fresh (unused) border/style slots, then
# clear them fully so no visible formatting remains and the dimension is unaffected.
$ws.Range("ZZ1000").Borders.Item(1).LineStyle = 1
$ws.Range("ZZ1000").Clear()
$ws.Range("ZZ1000").Borders.Item(1).LineStyle = 4
$ws.Range("ZZ1000").Clear()

# Header labels (text unchanged; rewritten so the string table is refreshed)
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "fedfundsrate"
$ws.Range("C1").Value = "fedfundsrate_shadow"

# Refreshed model output: updated fedfundsrate / fedfundsrate_shadow figures
$ws.Range("B2").Value = 9.6866666666665768
$ws.Range("C2").Value = 9.6866666666665768
$ws.Range("C3").Value = 10.556666666666569
$ws.Range("C4").Value = 11.389999999999878
$ws.Range("C5").Value = 9.2666666666666231
$ws.Range("C6").Value = 8.4766666666665991
$ws.Range("C7").Value = 7.9233333333332157
$ws.Range("C8").Value = 7.8999999999999293
$ws.Range("C9").Value = 8.1033333333332624
$ws.Range("C10").Value = 7.8266666666665596
$ws.Range("C11").Value = 6.9199999999998596
$ws.Range("C12").Value = 6.2066666666666048
$ws.Range("C13").Value = 6.266666666666576
$ws.Range("C14").Value = 6.2199999999999145
$ws.Range("C15").Value = 6.6499999999998893
$ws.Range("C16").Value = 6.8433333333332902
$ws.Range("C17").Value = 6.9166666666665266
$ws.Range("C18").Value = 6.6633333333332656
$ws.Range("C19").Value = 7.1566666666665668
$ws.Range("C20").Value = 7.9833333333332535
$ws.Range("C21").Value = 8.4699999999999331
$ws.Range("C22").Value = 9.4433333333332925
$ws.Range("C23").Value = 9.7266666666665724
$ws.Range("C24").Value = 9.0833333333331989
$ws.Range("C25").Value = 8.6133333333332516
$ws.Range("C26").Value = 8.2499999999999574
$ws.Range("C27").Value = 8.243333333333247
$ws.Range("C28").Value = 8.1599999999999451
$ws.Range("C29").Value = 7.7433333333332577
$ws.Range("C30").Value = 6.426666666666625
$ws.Range("C31").Value = 5.8633333333332205
$ws.Range("C32").Value = 5.6433333333332225
$ws.Range("C83").Value = 1.0095634079110516
$ws.Range("C84").Value = 1.4326601082768287
$ws.Range("C85").Value = 1.9491999475314881
$ws.Range("C86").Value = 2.4691438670969568
$ws.Range("C87").Value = 2.9424712147130272
$ws.Range("C88").Value = 3.4591678976461759
$ws.Range("C89").Value = 3.979224265190795
$ws.Range("C90").Value = 4.4559661970380526
$ws.Range("C91").Value = 4.9060541115272915
$ws.Range("C92").Value = 5.2461500897920788
$ws.Range("C93").Value = 5.2462507656800383
$ws.Range("C94").Value = 5.2563516530061527
$ws.Range("C95").Value = 5.2497838923445173
$ws.Range("C96").Value = 5.0732125765547531
$ws.Range("C97").Value = 4.4966362278598737
$ws.Range("C98").Value = 3.1767198540175512
$ws.Range("C99").Value = 2.0867960119344886
$ws.Range("C100").Value = 1.9401993364265957
$ws.Range("C101").Value = 0.50804939453934495
$ws.Range("C102").Value = 1.591392451496243
$ws.Range("C103").Value = 0.18153936519957892
$ws.Range("C104").Value = -0.73153365275446802
$ws.Range("C105").Value = -0.78973838638486393
$ws.Range("C106").Value = -0.61932171768126487
$ws.Range("C107").Value = -2.0881407141996711
$ws.Range("C108").Value = -1.916450710736417
$ws.Range("C109").Value = -2.560988738776182
$ws.Range("C110").Value = -2.0061879823781825
$ws.Range("C111").Value = -1.6079331322690571
$ws.Range("C112").Value = -2.8955959642015916
$ws.Range("C113").Value = -2.4625055418208697
$ws.Range("C114").Value = -3.3551770112617341
$ws.Range("C115").Value = -3.1330787178836306
$ws.Range("C116").Value = -2.6859582223774825
$ws.Range("C117").Value = -3.8952555099290964
$ws.Range("C118").Value = -2.3585523576021949
$ws.Range("C119").Value = -1.630369110900376
$ws.Range("C120").Value = -1.2308324022547446
$ws.Range("C121").Value = -1.3143042000176131
$ws.Range("C122").Value = -1.3608478644916544
$ws.Range("C123").Value = -1.2902938884155168
$ws.Range("C124").Value = -0.87651294332203866
$ws.Range("C125").Value = -0.41013593310598173
$ws.Range("C126").Value = 0.12832105241287195
$ws.Range("C127").Value = 0.016245456115782098
$ws.Range("C128").Value = 0.047772073577179164
$ws.Range("C147").Value = 7.7628117717102318
$ws.Range("C148").Value = -6.9697272932938485
$ws.Range("C149").Value = -4.5222797955217553
$ws.Range("C150").Value = -3.5240795397061486
$ws.Range("C151").Value = -3.1477835207301408
$ws.Range("C152").Value = -1.8116100516775258
$ws.Range("C153").Value = -0.64562495193367653

# Append latest observation (2022 Q1)
$ws.Range("A154").Value = 2022
$ws.Range("B154").Value = 0
$ws.Range("C154").Value = 0.32115971420862977
